# Regenerate merged AHB files
#
# The source/target column headers are renamed from the generic
# "_old" / "_new" suffixes to the concrete file-version tags
# "_FV2410" / "_FV2504", the header row is frozen, and the whole
# used range is wrapped in an Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the 20 header cells (columns A-J and L-U; K holds "diff"
#    and is left untouched).
$headers = @(
  "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410", "Segment ID_FV2410",
  "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410", "Bedingungsausdruck_FV2410", "Bedingung_FV2410",
  "diff",
  "Segmentname_FV2504", "Segmentgruppe_FV2504", "Segment_FV2504", "Datenelement_FV2504", "Segment ID_FV2504",
  "Code_FV2504", "Qualifier_FV2504", "Beschreibung_FV2504", "Bedingungsausdruck_FV2504", "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Wrap the whole used range in a native Excel table, picking up the
#    freshly renamed header text for the table's column names.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U72"), $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (split below row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
